$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Scalar fixes on existing rows
# ---------------------------------------------------------------------------
$ws.Range("Q59").Value  = 0    # was 2
$ws.Range("Q64").Value  = 0    # was 1
$ws.Range("O851").Value = 3    # was 0
$ws.Range("R853").Value = 0    # was blank
$ws.Range("R854").Value = 0    # was blank

# ---------------------------------------------------------------------------
# 2) Append the new weekly candles (rows 855-880), columns A-Q.
#    Column F (Adj Close) and column R (backup) stay blank for every new row,
#    exactly like the existing rows at the bottom of the sheet.
# ---------------------------------------------------------------------------
$dateFormat = $ws.Cells.Item(854, 1).NumberFormat   # reuse the existing date style

$newRows = @(
    @(855, 45474,64.97463376075062,68.2433028575731,64.12757329850771,67.32647705078125,$null,205957767,2024,7,1,0,0,0,27,0,0,0),
    @(856, 45481,67.66530155915338,69.5288329387574,64.67567121871676,68.12371063232422,$null,347775952,2024,7,8,0,0,0,28,0,0,0),
    @(857, 45488,68.5123629704095,71.75112921991978,66.41962322209983,66.69865417480469,$null,423534849,2024,7,15,0,0,0,29,0,0,0),
    @(858, 45495,66.56910982970626,69.30960449540036,61.28742616006817,66.68869018554688,$null,238521335,2024,7,22,0,0,0,30,0,0,0),
    @(859, 45502,67.11719696776022,67.43609055400719,63.55953730746447,64.26708221435547,$null,101230070,2024,7,29,0,0,0,31,0,0,0),
    @(860, 45509,61.82555989569462,64.21726582571256,60.09157252620486,61.75580215454102,$null,164768449,2024,8,5,0,0,0,32,0,0,0),
    @(861, 45516,61.75580330484236,63.57947597420266,61.53656251341199,62.62279510498047,$null,63844411,2024,8,12,0,0,0,33,0,0,0),
    @(862, 45523,62.88189250284754,65.3732512201,62.87192874064368,64.486328125,$null,78478577,2024,8,19,0,0,0,34,0,0,0),
    @(863, 45530,65.26468329240991,66.77179159814182,63.07888166957642,63.22859191894531,$null,107193214,2024,8,26,0,0,0,35,0,0,0),
    @(864, 45537,63.3783042866513,64.09692564810773,61.13261633948649,61.79135131835938,$null,91604935,2024,9,2,0,0,0,36,0,0,0),
    @(865, 45544,61.74144742628442,63.50805572129352,59.00669689427941,61.32225036621094,$null,111056979,2024,9,9,0,0,0,37,0,0,0),
    @(866, 45551,61.63165735105688,65.37447055051378,57.88884415159999,62.95911026000977,$null,216803058,2024,9,16,0,0,0,38,0,0,0),
    @(867, 45558,63.99712033641095,64.17677568103801,60.48386195587585,62.78943634033203,$null,70455867,2024,9,23,0,0,0,39,0,0,0),
    @(868, 45565,62.03089181925028,62.75949232642338,58.40784760313295,59.57560348510742,$null,59975413,2024,9,30,0,0,0,40,0,0,0),
    @(869, 45572,59.68539491566546,60.77330612629011,54.95448106937957,59.33606719970703,$null,79000643,2024,10,7,0,0,0,41,0,0,0),
    @(870, 45579,59.4857766236594,60.38405331996891,56.10227407987108,57.16024398803711,$null,52310454,2024,10,14,0,0,0,42,0,0,0),
    @(871, 45586,57.38980489690253,57.63932578775863,49.96406455568425,50.69266510009766,$null,82873460,2024,10,21,0,0,0,43,0,0,0),
    @(872, 45593,51.59093593159474,54.54526605319593,50.33335236790419,53.38748931884766,$null,98795757,2024,10,28,0,0,0,44,0,0,0),
    @(873, 45600,53.38748968739303,53.38748968739303,50.78249106527992,51.34141540527344,$null,62050607,2024,11,4,0,0,0,45,0,0,0),
    @(874, 45607,51.5,52.31999969482422,47.61000061035156,47.7599983215332,$null,57946354,2024,11,11,0,0,0,46,0,0,0),
    @(875, 45614,47.86999893188477,48.83000183105469,45.06000137329102,47.59999847412109,$null,69522748,2024,11,18,0,0,0,47,2,0,0),
    @(876, 45621,49.29999923706055,53.81999969482422,49.29999923706055,53.38999938964844,$null,121638174,2024,11,25,0,0,0,48,0,0,0),
    @(877, 45628,53.34000015258789,59.45999908447266,52.88999938964844,59,$null,98234712,2024,12,2,0,0,0,49,0,0,0),
    @(878, 45635,59.2400016784668,61.9900016784668,56.81000137329102,59.11000061035156,$null,211423180,2024,12,9,0,0,0,50,0,0,0),
    @(879, 45642,59.13999938964844,59.34000015258789,54.83000183105469,55.09000015258789,$null,60321321,2024,12,16,0,0,0,51,0,0,0),
    @(880, 45649,56,56.2400016784668,54.38000106811523,55.22999954223633,$null,40263358,2024,12,23,0,0,0,52,0,0,0)
)

foreach ($row in $newRows) {
    $r  = $row[0]
    $a  = $row[1]; $b  = $row[2];  $c  = $row[3];  $d  = $row[4];  $e  = $row[5]
    $f  = $row[6]; $g  = $row[7];  $h  = $row[8];  $i  = $row[9];  $j  = $row[10]
    $k  = $row[11]; $l = $row[12]; $m  = $row[13]; $n  = $row[14]; $o  = $row[15]
    $p  = $row[16]; $q = $row[17]

    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    if ($f -ne $null) { $ws.Cells.Item($r, 6).Value = $f }
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    # column R (18) is intentionally left blank for every new row
}
